$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3127422
$ws.Range("I40").Value = 8930114
$ws.Range("J40").Value = 2895.5386
$ws.Range("K40").Value = 8930114
$ws.Range("L40").Value = 2895.5386
$ws.Range("M40").Value = -8929939
$ws.Range("N40").Value = -3245.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12203650
$ws.Range("I32").Value = 14499700
$ws.Range("J32").Value = 16921.54
$ws.Range("K32").Value = 14499700
$ws.Range("L32").Value = 16921.54
$ws.Range("M32").Value = -14499413
$ws.Range("N32").Value = -17495.54

$ws.Range("H61").Value = 1597.1562
$ws.Range("I61").Value = 1412.5111
$ws.Range("J61").Value = 2034.4736
$ws.Range("K61").Value = 1412.5111
$ws.Range("L61").Value = 2034.4736
$ws.Range("M61").Value = -1200.5111
$ws.Range("N61").Value = -2458.4736

$ws.Range("H74").Value = 15626128
$ws.Range("I74").Value = 17242556
$ws.Range("J74").Value = 659.6667
$ws.Range("K74").Value = 17242556
$ws.Range("L74").Value = 659.6667
$ws.Range("M74").Value = -17241682
$ws.Range("N74").Value = -2407.6667

$ws.Range("H77").Value = 15626128
$ws.Range("I77").Value = 17242556
$ws.Range("J77").Value = 659.6667
$ws.Range("K77").Value = 86212780
$ws.Range("L77").Value = 3298.3335
$ws.Range("M77").Value = -86208412
$ws.Range("N77").Value = -12034.3335

$ws.Range("H132").Value = 22225614
$ws.Range("I132").Value = 27027996
$ws.Range("J132").Value = 14599.75
$ws.Range("K132").Value = 81083988
$ws.Range("L132").Value = 43799.25
$ws.Range("M132").Value = -81081458
$ws.Range("N132").Value = -48859.25

$ws.Range("H136").Value = 1597.1562
$ws.Range("I136").Value = 1412.5111
$ws.Range("J136").Value = 2034.4736
$ws.Range("K136").Value = 4237.5333
$ws.Range("L136").Value = 6103.4208
$ws.Range("M136").Value = -1687.5333
$ws.Range("N136").Value = -11203.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 971235.5600000001
$ws.Range("I86").Value = 2496.375
$ws.Range("J86").Value = 2908714
$ws.Range("K86").Value = 2496.375
$ws.Range("L86").Value = 2908714
$ws.Range("M86").Value = -1373.375
$ws.Range("N86").Value = -2910960

$ws.Range("H89").Value = 971235.5600000001
$ws.Range("I89").Value = 2496.375
$ws.Range("J89").Value = 2908714
$ws.Range("K89").Value = 12481.875
$ws.Range("L89").Value = 14543570
$ws.Range("M89").Value = -6865.875
$ws.Range("N89").Value = -14554802

$ws.Range("H134").Value = 1661196.9
$ws.Range("I134").Value = 969.73334
$ws.Range("J134").Value = 5057116
$ws.Range("K134").Value = 2909.20002
$ws.Range("L134").Value = 15171348
$ws.Range("M134").Value = -374.2000200000002
$ws.Range("N134").Value = -15176418

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 33334378
$ws.Range("I58").Value = 100001090
$ws.Range("J58").Value = 1022.65
$ws.Range("K58").Value = 100001090
$ws.Range("L58").Value = 1022.65
$ws.Range("M58").Value = -100000887
$ws.Range("N58").Value = -1428.65

$ws.Range("H132").Value = 18525834
$ws.Range("I132").Value = 1506.5
$ws.Range("J132").Value = 41681244
$ws.Range("K132").Value = 4519.5
$ws.Range("L132").Value = 125043732
$ws.Range("M132").Value = -1989.5
$ws.Range("N132").Value = -125048792

$ws.Range("H134").Value = 1687.3334
$ws.Range("I134").Value = 1352.8
$ws.Range("J134").Value = 1926.2858
$ws.Range("K134").Value = 4058.4
$ws.Range("L134").Value = 5778.857400000001
$ws.Range("M134").Value = -1523.4
$ws.Range("N134").Value = -10848.8574

$ws.Range("H136").Value = 33334378
$ws.Range("I136").Value = 100001090
$ws.Range("J136").Value = 1022.65
$ws.Range("K136").Value = 300003270
$ws.Range("L136").Value = 3067.95
$ws.Range("M136").Value = -300000720
$ws.Range("N136").Value = -8167.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 46671590
$ws.Range("I5").Value = 66667452
$ws.Range("J5").Value = 33341016
$ws.Range("K5").Value = 200002356
$ws.Range("L5").Value = 100023048
$ws.Range("M5").Value = -200002244
$ws.Range("N5").Value = -100023272

$ws.Range("H34").Value = 700.2857
$ws.Range("J34").Value = 2000
$ws.Range("L34").Value = 6000
$ws.Range("N34").Value = -6168

$ws.Range("H39").Value = 1150
$ws.Range("I39").Value = 512.5
$ws.Range("K39").Value = 1537.5
$ws.Range("M39").Value = -1243.5

$ws.Range("H55").Value = 1665.3846
$ws.Range("J55").Value = 1804.5454
$ws.Range("L55").Value = 5413.6362
$ws.Range("N55").Value = -5767.6362

$ws.Range("H131").Value = 707.8461
$ws.Range("J131").Value = 780.8767
$ws.Range("L131").Value = 2342.6301
$ws.Range("N131").Value = -12422.6301

$ws.Range("H135").Value = 46671590
$ws.Range("I135").Value = 66667452
$ws.Range("J135").Value = 33341016
$ws.Range("K135").Value = 600007068
$ws.Range("L135").Value = 300069144
$ws.Range("M135").Value = -600004533
$ws.Range("N135").Value = -300074214

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4589.857
$ws.Range("I132").Value = 832.5909
$ws.Range("J132").Value = 18366.5
$ws.Range("K132").Value = 2497.7727
$ws.Range("L132").Value = 55099.5
$ws.Range("M132").Value = 32.22730000000001
$ws.Range("N132").Value = -60159.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1504.3158
$ws.Range("I68").Value = 1513
$ws.Range("J68").Value = 1480
$ws.Range("K68").Value = 1513
$ws.Range("L68").Value = 1480
$ws.Range("M68").Value = -764
$ws.Range("N68").Value = -2978

$ws.Range("H71").Value = 1504.3158
$ws.Range("I71").Value = 1513
$ws.Range("J71").Value = 1480
$ws.Range("K71").Value = 7565
$ws.Range("L71").Value = 7400
$ws.Range("M71").Value = -3821
$ws.Range("N71").Value = -14888

$ws.Range("H136").Value = 3878.3057
$ws.Range("I136").Value = 4081.5386
$ws.Range("J136").Value = 3349.9
$ws.Range("K136").Value = 12244.6158
$ws.Range("L136").Value = 10049.7
$ws.Range("M136").Value = -9694.6158
$ws.Range("N136").Value = -15149.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2975
$ws.Range("I62").Value = 2970
$ws.Range("K62").Value = 2970
$ws.Range("M62").Value = -2346

$ws.Range("H65").Value = 2975
$ws.Range("I65").Value = 2970
$ws.Range("K65").Value = 14850
$ws.Range("M65").Value = -11730

$ws.Range("H123").Value = 20102.285
$ws.Range("J123").Value = 20102.285
$ws.Range("L123").Value = 20102.285
$ws.Range("N123").Value = -29902.285

$ws.Range("H132").Value = 20852.21
$ws.Range("I132").Value = 30196.943
$ws.Range("J132").Value = 5985.591
$ws.Range("K132").Value = 90590.829
$ws.Range("L132").Value = 17956.773
$ws.Range("M132").Value = -88060.829
$ws.Range("N132").Value = -23016.773

$ws.Range("H136").Value = 3493.9792
$ws.Range("I136").Value = 5331.7827
$ws.Range("J136").Value = 1803.2
$ws.Range("K136").Value = 15995.3481
$ws.Range("L136").Value = 5409.6
$ws.Range("M136").Value = -13445.3481
$ws.Range("N136").Value = -10509.6
